# Daily attendance processing - 2025-10-21 08:52:16
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) "Recorded By" (column G) value re-ordering fixes ---
# The recorder list for these sessions was re-synced from the source system,
# which re-orders the comma-joined list of who recorded attendance.
$ws.Range("G2").Value = 'system, backup@backdoor.com, System'
$ws.Range("G3").Value = 'dnasr281@gmail.com, System'
$ws.Range("G6").Value = 'dnasr281@gmail.com, System'
$ws.Range("G10").Value = 'dnasr281@gmail.com, System'
$ws.Range("G12").Value = 'dnasr281@gmail.com, System'
$ws.Range("G13").Value = 'dnasr281@gmail.com, System'
$ws.Range("G14").Value = 'dnasr281@gmail.com, System'
$ws.Range("G15").Value = 'dnasr281@gmail.com, System'
$ws.Range("G18").Value = 'dnasr281@gmail.com, System'
$ws.Range("G19").Value = 'dnasr281@gmail.com, System'
$ws.Range("G20").Value = 'dnasr281@gmail.com, System'
$ws.Range("G21").Value = 'dnasr281@gmail.com, System'
$ws.Range("G22").Value = 'dnasr281@gmail.com, System'
$ws.Range("G29").Value = 'system, backup@backdoor.com, System'
$ws.Range("G30").Value = 'dnasr281@gmail.com, System'
$ws.Range("G33").Value = 'dnasr281@gmail.com, System'
$ws.Range("G37").Value = 'dnasr281@gmail.com, System'
$ws.Range("G39").Value = 'dnasr281@gmail.com, System'
$ws.Range("G40").Value = 'dnasr281@gmail.com, System'
$ws.Range("G41").Value = 'dnasr281@gmail.com, System'
$ws.Range("G42").Value = 'dnasr281@gmail.com, System'
$ws.Range("G45").Value = 'dnasr281@gmail.com, System'
$ws.Range("G46").Value = 'dnasr281@gmail.com, System'
$ws.Range("G47").Value = 'dnasr281@gmail.com, System'
$ws.Range("G48").Value = 'dnasr281@gmail.com, System'
$ws.Range("G49").Value = 'dnasr281@gmail.com, System'
$ws.Range("G56").Value = 'system, backup@backdoor.com, System'
$ws.Range("G57").Value = 'dnasr281@gmail.com, System'
$ws.Range("G60").Value = 'dnasr281@gmail.com, System'
$ws.Range("G64").Value = 'dnasr281@gmail.com, System'
$ws.Range("G66").Value = 'dnasr281@gmail.com, System'
$ws.Range("G67").Value = 'dnasr281@gmail.com, System'
$ws.Range("G68").Value = 'dnasr281@gmail.com, System'
$ws.Range("G69").Value = 'dnasr281@gmail.com, System'
$ws.Range("G72").Value = 'dnasr281@gmail.com, System'
$ws.Range("G73").Value = 'dnasr281@gmail.com, System'
$ws.Range("G74").Value = 'dnasr281@gmail.com, System'
$ws.Range("G75").Value = 'dnasr281@gmail.com, System'
$ws.Range("G76").Value = 'dnasr281@gmail.com, System'
$ws.Range("G86").Value = 'dnasr281@gmail.com, System'
$ws.Range("G87").Value = 'dnasr281@gmail.com, System'
$ws.Range("G88").Value = 'dnasr281@gmail.com, System'
$ws.Range("G89").Value = 'dnasr281@gmail.com, System'
$ws.Range("G93").Value = 'dnasr281@gmail.com, System'
$ws.Range("G95").Value = 'dnasr281@gmail.com, System'
$ws.Range("G99").Value = 'dnasr281@gmail.com, System'
$ws.Range("G112").Value = 'dnasr281@gmail.com, System'
$ws.Range("G113").Value = 'dnasr281@gmail.com, System'
$ws.Range("G114").Value = 'dnasr281@gmail.com, System'
$ws.Range("G115").Value = 'dnasr281@gmail.com, System'
$ws.Range("G119").Value = 'dnasr281@gmail.com, System'
$ws.Range("G121").Value = 'dnasr281@gmail.com, System'
$ws.Range("G125").Value = 'dnasr281@gmail.com, System'
$ws.Range("G138").Value = 'dnasr281@gmail.com, System'
$ws.Range("G139").Value = 'dnasr281@gmail.com, System'
$ws.Range("G140").Value = 'dnasr281@gmail.com, System'
$ws.Range("G141").Value = 'dnasr281@gmail.com, System'
$ws.Range("G145").Value = 'dnasr281@gmail.com, System'
$ws.Range("G147").Value = 'dnasr281@gmail.com, System'
$ws.Range("G151").Value = 'dnasr281@gmail.com, System'

# --- 2) Newly-recorded sessions (row 100 / B2D, row 126 / B2E, row 152 / B2F) ---
# These three sessions on 21/10/2025 were "Not Recorded" and are now recorded;
# copy the formatting (fill/font) used by already-recorded rows so the pink
# "not recorded" highlight style is replaced with the normal row style, and
# fill in recorder / attendance counts / status.
$fmtSource = $ws.Range("A99:I99")
$fmtSource.Copy()
$ws.Range("A100:I100").PasteSpecial(-4122)
$ws.Range("A126:I126").PasteSpecial(-4122)
$ws.Range("A152:I152").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("G100").Value = "dnasr281@gmail.com"
$ws.Range("H100").Value = "40/56"
$ws.Range("I100").Value = "Recorded"

$ws.Range("G126").Value = "dnasr281@gmail.com"
$ws.Range("H126").Value = "31/55"
$ws.Range("I126").Value = "Recorded"

$ws.Range("G152").Value = "dnasr281@gmail.com"
$ws.Range("H152").Value = "45/57"
$ws.Range("I152").Value = "Recorded"

# --- 3) Class Statistics summary block (K3:L10) ---
$ws.Range("L6").Value = 117
$ws.Range("L7").Value = 0
$ws.Range("L9").Value = "73.6%"
$ws.Range("L10").Value = "72.2%"

# --- 4) Group Statistics block (rows 18-20: B2D/B2E/B2F) recalculated with the ---
# --- newly recorded 21/10/2025 sessions folded in ---
$ws.Range("O18").Value = 18
$ws.Range("P18").Value = 0
$ws.Range("R18").Value = "69.2%"
$ws.Range("S18").Value = "78.0%"

$ws.Range("O19").Value = 18
$ws.Range("P19").Value = 0
$ws.Range("R19").Value = "69.2%"
$ws.Range("S19").Value = "75.9%"

$ws.Range("O20").Value = 18
$ws.Range("P20").Value = 0
$ws.Range("R20").Value = "69.2%"
$ws.Range("S20").Value = "82.9%"

# --- 5) Column I ("Students") width tightened ---
$ws.Columns("I").ColumnWidth = 10

Write-Host "Daily attendance processing complete."
